# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

# Mapping of row number -> new value for column F.
$updates = @{
    2  = 3430
    5  = 1768
    6  = 1672
    8  = 391
    13 = 240
    14 = 15
    15 = 68
    21 = 67
    22 = 147
    24 = 425
    25 = 311
    27 = 52
    28 = 22
    29 = 33
    30 = 627
    31 = 2426
    32 = 17
    33 = 54
    35 = 699
    38 = 244
    39 = 367
    41 = 564
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
